# Daily attendance processing
# Normalize the "Recorded By" (column G) entries so that "System" always
# appears first in the comma-separated list of recorders, preserving the
# relative order of the remaining entries.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "backup@backdoor.com, system, System" -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -eq "System") {
            $rest = $parts[0..($parts.Count - 2)]
            $newParts = @("System") + $rest
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
